$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set this first so "ADC1CH2" lands in the shared-strings table ahead of
# "pre fix" / "post fix" (matches the target's shared-string ordering).
$ws.Range("J18").Value = "ADC1CH2"

$ws.Range("A1").Value = "pre fix"
$ws.Range("G1").Value = "post fix"

$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 1
$ws.Range("K9").Value = "PA0"

$ws.Range("I11").Value = 4
$ws.Range("K11").Value = "PB6"

$ws.Range("I12").Value = 4
$ws.Range("K12").Value = "PB7"

$ws.Range("I18").Value = "x"
$ws.Range("K18").Value = "PA1"

$ws.Range("K22").Value = "PA10"
$ws.Range("K23").Value = "PB3"
$ws.Range("K24").Value = "PB4"

$ws.Range("H13").Select()
